$d = $word.ActiveDocument

function Escape-Xml([string]$s) {
    $s = $s -replace '&', '&amp;'
    $s = $s -replace '<', '&lt;'
    $s = $s -replace '>', '&gt;'
    return $s
}

# ---------------------------------------------------------------------------
# Edit 1: "Week beginning Monday April 12th" cell - append a new run with a
# sentence about the Zoom meeting to the existing paragraph.
# ---------------------------------------------------------------------------
$findRng1 = $d.Content
$found1 = $findRng1.Find.Execute(
    "Added Clara as collaborator.  Assisted John with Git workflow issue about branches he was having.  Reviewed PRs from Andrea.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found1) {
    throw "Could not locate the April 12th update paragraph."
}

$existingText1 = Escape-Xml $findRng1.Text
$xml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' +
    '<w:p>' +
    '<w:r><w:t>' + $existingText1 + '</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> Held Zoom meeting with group to organize thoughts and plans for the week.</w:t></w:r>' +
    '</w:p>' +
    '</w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$findRng1.InsertXML($xml1)

# ---------------------------------------------------------------------------
# Edit 2: "Week beginning Monday April 19th" cell - the first blank paragraph
# right after the heading becomes a paragraph with two runs describing the
# week's work.
# ---------------------------------------------------------------------------
$findRng2 = $d.Content
$found2 = $findRng2.Find.Execute("April 19th", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found2) {
    throw "Could not locate the April 19th heading."
}

# Figure out which document paragraph holds the "April 19th" heading, then
# grab the very next paragraph (the blank one to be filled in).
$headingIndex = 0
$i = 1
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Start -le $findRng2.Start -and $p.Range.End -ge $findRng2.End) {
        $headingIndex = $i
    }
    $i++
}

if ($headingIndex -eq 0) {
    throw "Could not map the April 19th heading to a paragraph index."
}

$targetPara = $d.Paragraphs.Item($headingIndex + 1)
$targetRng = $targetPara.Range

if ($targetRng.Text -ne "" -and $targetRng.Text -ne "`r") {
    throw "Expected an empty paragraph after the April 19th heading, found: $($targetRng.Text)"
}

$weekText = "Deployed LMNOP to GCP. Helped John realize he was missing a requirements.txt installation when testing Andrea" +
    [char]0x2019 + "s pagination code.  Brainstormed what would go into a user" + [char]0x2019 + "s profile page with Andrea"

$xml2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' +
    '<w:p>' +
    '<w:r><w:t>' + $weekText + '</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> and what would be visible to the public vs. account owner.</w:t></w:r>' +
    '</w:p>' +
    '</w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$targetRng.InsertXML($xml2)

Write-Output "done"
